$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice")

# Write discount explanation text (replaces the numeric unit price 1 in E18)
$ws.Range("E18").Value = "This client doesn't benefit from any discount"

# Write signature text to A31 (replaces placeholder "name, email address")
$ws.Range("A31").Value = "RPA Dev, developer.rpa@mail.com"
